$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019321692680165
$ws.Range("C2").Value = 0.0479486225917185
$ws.Range("D2").Value = 0.0991888149321003
$ws.Range("E2").Value = 0.05874338308593607
$ws.Range("F2").Value = 1.92141472027599
$ws.Range("I2").Value = 1.544883511641942
$ws.Range("K2").Value = 0.8148757401611704
$ws.Range("L2").Value = 0.2467597928773984
$ws.Range("N2").Value = 2.627067627999558
$ws.Range("B3").Value = 0.9749936956407907
$ws.Range("C3").Value = 0.04172204689865566
$ws.Range("D3").Value = 0.09938035963789105
$ws.Range("E3").Value = 0.05837070214493423
$ws.Range("F3").Value = 1.907890773193515
$ws.Range("I3").Value = 1.542902292285319
$ws.Range("K3").Value = 0.7674244833306147
$ws.Range("L3").Value = 0.239294441631003
$ws.Range("N3").Value = 2.641760148719186
$ws.Range("B4").Value = 0.9483559185673585
$ws.Range("C4").Value = 0.03789546538128263
$ws.Range("D4").Value = 0.09949987815980688
$ws.Range("E4").Value = 0.05816767181454274
$ws.Range("F4").Value = 1.900532097282081
$ws.Range("I4").Value = 1.542342171065279
$ws.Range("K4").Value = 0.7387516762821065
$ws.Range("L4").Value = 0.2348535731390058
$ws.Range("N4").Value = 2.651500712145435
$ws.Range("B5").Value = 0.9376467304189759
$ws.Range("C5").Value = 0.03633515587704039
$ws.Range("D5").Value = 0.09954906736362723
$ws.Range("E5").Value = 0.05809143077961565
$ws.Range("F5").Value = 1.897770825965011
$ws.Range("I5").Value = 1.542278880449466
$ws.Range("K5").Value = 0.7271835219089553
$ws.Range("L5").Value = 0.2330797952925678
$ws.Range("N5").Value = 2.655650810055697
$ws.Range("B6").Value = 0.9358772951912044
$ws.Range("C6").Value = 0.03607600812971157
$ws.Range("D6").Value = 0.09955726458235148
$ws.Range("E6").Value = 0.05807916357835019
$ws.Range("F6").Value = 1.897326654564054
$ws.Range("I6").Value = 1.542278332122656
$ws.Range("K6").Value = 0.7252696587093226
$ws.Range("L6").Value = 0.2327874293393819
$ws.Range("N6").Value = 2.656350843001192
$ws.Range("B7").Value = 0.9482108997030139
$ws.Range("C7").Value = 0.03787442644021155
$ws.Range("D7").Value = 0.09950053957581595
$ws.Range("E7").Value = 0.05816661728979788
$ws.Range("F7").Value = 1.900493896660748
$ws.Range("I7").Value = 1.542340649690843
$ws.Range("K7").Value = 0.7385951934766979
$ws.Range("L7").Value = 0.2348295059659193
$ws.Range("N7").Value = 2.65155595020002
$ws.Range("B8").Value = 1.003917122852499
$ws.Range("C8").Value = 0.04580236202707511
$ws.Range("D8").Value = 0.09925446586268194
$ws.Range("E8").Value = 0.05860953331864849
$ws.Range("F8").Value = 1.916555381448362
$ws.Range("I8").Value = 1.544064106237634
$ws.Range("K8").Value = 0.798418521257048
$ws.Range("L8").Value = 0.2441560694356184
$ws.Range("N8").Value = 2.631984238538344
$ws.Range("B9").Value = 1.117758773852671
$ws.Range("C9").Value = 0.06132614393251856
$ws.Range("D9").Value = 0.09878686389830449
$ws.Range("E9").Value = 0.0596825176538438
$ws.Range("F9").Value = 1.955563800163304
$ws.Range("I9").Value = 1.552657361778245
$ws.Range("K9").Value = 0.9194104010799151
$ws.Range("L9").Value = 0.2635813068848165
$ws.Range("N9").Value = 2.599317529893085
$ws.Range("B10").Value = 1.204216383256778
$ws.Range("C10").Value = 0.07272466747247108
$ws.Range("D10").Value = 0.09845214126901425
$ws.Range("E10").Value = 0.06059527786079144
$ws.Range("F10").Value = 1.988826758970973
$ws.Range("I10").Value = 1.562159834447073
$ws.Range("K10").Value = 1.010568793196654
$ws.Range("L10").Value = 0.2785503289682794
$ws.Range("N10").Value = 2.57880877176747
$ws.Range("B11").Value = 1.244163574549759
$ws.Range("C11").Value = 0.0779102765572901
$ws.Range("D11").Value = 0.09830172423475858
$ws.Range("E11").Value = 0.06103751007249869
$ws.Range("F11").Value = 2.004964152080149
$ws.Range("I11").Value = 1.567177751167932
$ws.Range("K11").Value = 1.052536748941037
$ws.Range("L11").Value = 0.2855127060836082
$ws.Range("N11").Value = 2.570238861048011
$ws.Range("B12").Value = 1.259379349813287
$ws.Range("C12").Value = 0.07987409019077063
$ws.Range("D12").Value = 0.09824502730949192
$ws.Range("E12").Value = 0.06120885095187845
$ws.Range("F12").Value = 2.011219958137588
$ws.Range("I12").Value = 1.569178034118039
$ws.Range("K12").Value = 1.06850099525596
$ws.Range("L12").Value = 0.2881712189370091
$ws.Range("N12").Value = 2.567103070455417
$ws.Range("B13").Value = 1.256098419745456
$ws.Range("C13").Value = 0.07945114040069257
$ws.Range("D13").Value = 0.09825722637188328
$ws.Range("E13").Value = 0.06117177730625656
$ws.Range("F13").Value = 2.009866208305965
$ws.Range("I13").Value = 1.568742782876583
$ws.Range("K13").Value = 1.065059610772835
$ws.Range("L13").Value = 0.2875976812304515
$ws.Range("N13").Value = 2.56777354783415
$ws.Range("B14").Value = 1.24541360965884
$ws.Range("C14").Value = 0.0780718374790581
$ws.Range("D14").Value = 0.0982970544967614
$ws.Range("E14").Value = 0.06105152872763853
$ws.Range("F14").Value = 2.005475914882538
$ws.Range("I14").Value = 1.567340308556268
$ws.Range("K14").Value = 1.053848695915917
$ws.Range("L14").Value = 0.2857309822372258
$ws.Range("N14").Value = 2.569978682908356
$ws.Range("B15").Value = 1.238880395424587
$ws.Range("C15").Value = 0.07722699481864481
$ws.Range("D15").Value = 0.09832148449686962
$ws.Range("E15").Value = 0.0609783777038615
$ws.Range("F15").Value = 2.002805615278575
$ws.Range("I15").Value = 1.56649429356473
$ws.Range("K15").Value = 1.046991051157477
$ws.Range("L15").Value = 0.2845904425619921
$ws.Range("N15").Value = 2.571343651710762
$ws.Range("B16").Value = 1.201618150276602
$ws.Range("C16").Value = 0.07238579018837754
$ws.Range("D16").Value = 0.09846200831768925
$ws.Range("E16").Value = 0.06056691989040175
$ws.Range("F16").Value = 1.987792407402935
$ws.Range("I16").Value = 1.561845902663293
$ws.Range("K16").Value = 1.007836145608536
$ws.Range("L16").Value = 0.2780984000527837
$ws.Range("N16").Value = 2.579384139516222
$ws.Range("B17").Value = 1.178916967219266
$ws.Range("C17").Value = 0.06941602052307871
$ws.Range("D17").Value = 0.09854868655772542
$ws.Range("E17").Value = 0.06032141730407403
$ws.Range("F17").Value = 1.978840123445622
$ws.Range("I17").Value = 1.559172419105124
$ws.Range("K17").Value = 0.9839438650504917
$ws.Range("L17").Value = 0.2741549201518438
$ws.Range("N17").Value = 2.584511428736917
$ws.Range("B18").Value = 1.165917946412264
$ws.Range("C18").Value = 0.06770791675278076
$ws.Range("D18").Value = 0.09859871589973856
$ws.Range("E18").Value = 0.06018275373365256
$ws.Range("F18").Value = 1.973785670279469
$ws.Range("I18").Value = 1.557700128545797
$ws.Range("K18").Value = 0.9702486730898841
$ws.Range("L18").Value = 0.2719011243273144
$ws.Range("N18").Value = 2.587531996320962
$ws.Range("B19").Value = 1.161526683367981
$ws.Range("C19").Value = 0.06712958463407404
$ws.Range("D19").Value = 0.0986156850166573
$ws.Range("E19").Value = 0.06013624162191888
$ws.Range("F19").Value = 1.97209056950868
$ws.Range("I19").Value = 1.557212869150085
$ws.Range("K19").Value = 0.9656197897278957
$ws.Range("L19").Value = 0.2711404984242449
$ws.Range("N19").Value = 2.588566979993359
$ws.Range("B20").Value = 1.181327534668696
$ws.Range("C20").Value = 0.06973215394860688
$ws.Range("D20").Value = 0.09853944149450378
$ws.Range("E20").Value = 0.060347288333503
$ws.Range("F20").Value = 1.979783310004962
$ws.Range("I20").Value = 1.559450243596984
$ws.Range("K20").Value = 0.9864823727289433
$ws.Range("L20").Value = 0.2745732208760501
$ws.Range("N20").Value = 2.58395821990851
$ws.Range("B21").Value = 1.248549592554355
$ws.Range("C21").Value = 0.07847696782175717
$ws.Range("D21").Value = 0.0982853489057991
$ws.Range("E21").Value = 0.06108674348318033
$ws.Range("F21").Value = 2.006761515013963
$ws.Range("I21").Value = 1.567749531270636
$ws.Range("K21").Value = 1.057139661744458
$ws.Range("L21").Value = 0.2862786793793788
$ws.Range("N21").Value = 2.569328009143135
$ws.Range("B22").Value = 1.292999712348433
$ws.Range("C22").Value = 0.08419304116171134
$ws.Range("D22").Value = 0.09812081477794798
$ws.Range("E22").Value = 0.06159261578525843
$ws.Range("F22").Value = 2.025238200155101
$ws.Range("I22").Value = 1.573757161003243
$ws.Range("K22").Value = 1.103737358089688
$ws.Range("L22").Value = 0.2940571869794582
$ws.Range("N22").Value = 2.56040437376258
$ws.Range("B23").Value = 1.269228611475455
$ws.Range("C23").Value = 0.08114216134545416
$ws.Range("D23").Value = 0.09820849077439853
$ws.Range("E23").Value = 0.06132055710232365
$ws.Range("F23").Value = 2.015299447268021
$ws.Range("I23").Value = 1.570497333083694
$ws.Range("K23").Value = 1.078828924985572
$ws.Range("L23").Value = 0.2898939016380098
$ws.Range("N23").Value = 2.565108631819228
$ws.Range("B24").Value = 1.180237554544419
$ws.Range("C24").Value = 0.06958923235227132
$ws.Range("D24").Value = 0.09854362057342669
$ws.Range("E24").Value = 0.06033558431882469
$ws.Range("F24").Value = 1.979356607754454
$ws.Range("I24").Value = 1.559324437510398
$ws.Range("K24").Value = 0.9853345863319873
$ws.Range("L24").Value = 0.2743840654979834
$ws.Range("N24").Value = 2.584208098779115
$ws.Range("B25").Value = 1.086467629072359
$ws.Range("C25").Value = 0.05712826127430048
$ws.Range("D25").Value = 0.09891179283131457
$ws.Range("E25").Value = 0.05937038119141391
$ws.Range("F25").Value = 1.94420432214568
$ws.Range("I25").Value = 1.54977341774029
$ws.Range("K25").Value = 0.8862827057751588
$ws.Range("L25").Value = 0.2582041647421676
$ws.Range("N25").Value = 2.607542227327897
